$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# (e.g. "215.55") are preserved verbatim instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.814.27"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "1.628.70"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "215.55"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "0.5062"
$ws.Range("E6").Value = "  -0.63%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "0.06457"
$ws.Range("E8").Value = "  +1.52%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "0.2579"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("E10").Value = "  -2.56%  "
$ws.Range("D11").Value = "0.07805"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").Value = "4.257"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").Value = "1.628.00"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").Value = "1.853.76"
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("D15").Value = "0.5579"
$ws.Range("E15").Value = "  +1.85%  "
$ws.Range("D16").Value = "63.23"
$ws.Range("E16").Value = "  -1.69%  "
$ws.Range("D17").Value = "0.0₅7529"
$ws.Range("E17").Value = "  -2.96%  "
$ws.Range("D18").Value = "25.822.05"
$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").Value = "193.37"
$ws.Range("E20").Value = "  -1.55%  "
$ws.Range("E21").Value = "  -3.00%  "
$ws.Range("D22").Value = "9.805"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("D23").Value = "6.004"
$ws.Range("E23").Value = "  -1.38%  "
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("D25").Value = "1.816"
$ws.Range("E25").Value = "  -3.95%  "
$ws.Range("D26").Value = "140.21"
$ws.Range("E26").Value = "  -2.46%  "
$ws.Range("D27").Value = "0.1265"
$ws.Range("E27").Value = "  +2.36%  "
$ws.Range("D28").Value = "6.715"
$ws.Range("E28").Value = "  -2.30%  "
$ws.Range("D29").Value = "15.39"
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("D30").Value = "1.238"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "0.04853"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("D33").Value = "3.180"
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("D34").Value = "1.554"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("D36").Value = "0.8938"
$ws.Range("E36").Value = "  -2.36%  "
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "2.570"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.132.34"
$ws.Range("E38").Value = "  +3.85%  "
$ws.Range("D39").Value = "0.5466"
$ws.Range("E39").Value = "  -1.53%  "
$ws.Range("D40").Value = "0.01556"
$ws.Range("E40").Value = "  -1.05%  "
$ws.Range("D41").Value = "0.9994"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").Value = "5.563"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("D43").Value = "0.7957"
$ws.Range("E43").Value = "  -1.26%  "
$ws.Range("D44").Value = "97.27"
$ws.Range("E44").Value = "  -1.86%  "
$ws.Range("D45").Value = "1.782.01"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("D46").Value = "0.0₈111"
$ws.Range("E46").Value = "  -7.91%  "
$ws.Range("D47").Value = "0.4440"
$ws.Range("E47").Value = "  -2.13%  "
$ws.Range("D48").Value = "54.97"
$ws.Range("E48").Value = "  -0.97%  "
$ws.Range("D49").Value = "0.05059"
$ws.Range("D50").Value = "7.618"
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").Value = "1.002"
$ws.Range("E51").Value = "  -0.30%  "

# Restore the default (unformatted) style on column D so no stray
# number-format style lingers on these cells.
$ws.Range("D2:D51").Style = "Normal"

